$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A to make room for the "Tests" labels,
# shifting the existing Username/Password data from A:B to B:C.
$ws.Range("A1").EntireColumn.Insert()

# New column A: test case labels
$ws.Range("A1").Value = "Tests"
$ws.Range("A2").Value = "Invalid 1"
$ws.Range("A3").Value = "Invalid 2"
$ws.Range("A4").Value = "Invalid 3"
$ws.Range("A5").Value = "Valid 1"

# Rename the sheet
$ws.Name = "Login"

# Column widths to match the authored layout (A: custom width, B: best-fit)
$ws.Columns.Item(1).ColumnWidth = 12.8333333333
$ws.Columns.Item(2).ColumnWidth = 10

# Selection as left by the author
$ws.Range("D3").Select()

# Window size as recorded in the saved workbook
$excel.ActiveWindow.Width = 16830
$excel.ActiveWindow.Height = 2520
